$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="27.048.75"'
$ws.Range("E2").Formula = '="  +2.39%  "'
$ws.Range("D3").Formula = '="1.820.54"'
$ws.Range("E3").Formula = '="  +2.85%  "'
$ws.Range("D4").Formula = '="1.009"'
$ws.Range("E4").Formula = '="  +0.65%  "'
$ws.Range("D5").Formula = '="314.90"'
$ws.Range("E5").Formula = '="  +2.81%  "'
$ws.Range("D6").Formula = '="1.009"'
$ws.Range("E6").Formula = '="  +0.69%  "'
$ws.Range("D7").Formula = '="0.4305"'
$ws.Range("E7").Formula = '="  +0.56%  "'
$ws.Range("D8").Formula = '="0.3686"'
$ws.Range("E8").Formula = '="  +1.03%  "'
$ws.Range("D9").Formula = '="0.07270"'
$ws.Range("E9").Formula = '="  +1.31%  "'
$ws.Range("D10").Formula = '="2.186.71"'
$ws.Range("E10").Formula = '="  +21.71%  "'
$ws.Range("D11").Formula = '="0.8678"'
$ws.Range("E11").Formula = '="  +2.34%  "'
$ws.Range("D12").Formula = '="21.30"'
$ws.Range("E12").Formula = '="  +4.50%  "'
$ws.Range("D13").Formula = '="5.428"'
$ws.Range("E13").Formula = '="  +3.61%  "'
$ws.Range("D14").Formula = '="6.627"'
$ws.Range("E14").Formula = '="  +2.91%  "'
$ws.Range("D15").Formula = '="0.06969"'
$ws.Range("E15").Formula = '="  +0.92%  "'
$ws.Range("D16").Formula = '="81.43"'
$ws.Range("E16").Formula = '="  +2.76%  "'
$ws.Range("D17").Formula = '="1.015"'
$ws.Range("E17").Formula = '="  +0.97%  "'
$ws.Range("D18").Formula = '="0.000008943"'
$ws.Range("E18").Formula = '="  +3.52%  "'
$ws.Range("D19").Formula = '="1.008"'
$ws.Range("E19").Formula = '="  +0.57%  "'
$ws.Range("D20").Formula = '="15.31"'
$ws.Range("E20").Formula = '="  +2.03%  "'
$ws.Range("D21").Formula = '="27.083.71"'
$ws.Range("E21").Formula = '="  +2.49%  "'
$ws.Range("D22").Formula = '="5.180"'
$ws.Range("E22").Formula = '="  +1.27%  "'
$ws.Range("D23").Formula = '="2.401.33"'
$ws.Range("E23").Formula = '="  +20.19%  "'
$ws.Range("D24").Formula = '="11.01"'
$ws.Range("E24").Formula = '="  -1.30%  "'
$ws.Range("D25").Formula = '="154.48"'
$ws.Range("E25").Formula = '="  +1.61%  "'
$ws.Range("D26").Formula = '="1.882"'
$ws.Range("E26").Formula = '="  +1.14%  "'
$ws.Range("D27").Formula = '="18.35"'
$ws.Range("E27").Formula = '="  +1.67%  "'
$ws.Range("D28").Formula = '="5.252"'
$ws.Range("E28").Formula = '="  +3.56%  "'
$ws.Range("D29").Formula = '="1.909"'
$ws.Range("E29").Formula = '="  +8.98%  "'
$ws.Range("D30").Formula = '="114.76"'
$ws.Range("E30").Formula = '="  +0.35%  "'
$ws.Range("D31").Formula = '="0.08971"'
$ws.Range("E31").Formula = '="  +0.07%  "'
$ws.Range("D32").Formula = '="1.193"'
$ws.Range("E32").Formula = '="  +7.12%  "'
$ws.Range("D33").Formula = '="0.7553"'
$ws.Range("E33").Formula = '="  +4.29%  "'
$ws.Range("D34").Formula = '="4.436"'
$ws.Range("E34").Formula = '="  +2.49%  "'
$ws.Range("D35").Formula = '="2.816"'
$ws.Range("E35").Formula = '="  +2.66%  "'
$ws.Range("D36").Formula = '="1.009"'
$ws.Range("E36").Formula = '="  +0.72%  "'
$ws.Range("D37").Formula = '="1.136"'
$ws.Range("E37").Formula = '="  +5.24%  "'
$ws.Range("D38").Formula = '="0.05230"'
$ws.Range("E38").Formula = '="  +1.25%  "'
$ws.Range("D39").Formula = '="0.01933"'
$ws.Range("E39").Formula = '="  +2.45%  "'
$ws.Range("D40").Formula = '="0.5123"'
$ws.Range("E40").Formula = '="  +4.25%  "'
$ws.Range("D41").Formula = '="2.755"'
$ws.Range("E41").Formula = '="  +6.69%  "'
$ws.Range("D42").Formula = '="0.1656"'
$ws.Range("E42").Formula = '="  +2.64%  "'
$ws.Range("D43").Formula = '="6.511"'
$ws.Range("E43").Formula = '="  +3.48%  "'
$ws.Range("D44").Formula = '="8.361"'
$ws.Range("E44").Formula = '="  +4.23%  "'
$ws.Range("D45").Formula = '="107.15"'
$ws.Range("E45").Formula = '="  +2.31%  "'
$ws.Range("D46").Formula = '="10.41"'
$ws.Range("E46").Formula = '="  +2.55%  "'
$ws.Range("D47").Formula = '="1.010"'
$ws.Range("E47").Formula = '="  +0.87%  "'
$ws.Range("D48").Formula = '="0.4595"'
$ws.Range("E48").Formula = '="  +2.83%  "'
$ws.Range("D49").Formula = '="1.653"'
$ws.Range("E49").Formula = '="  +3.70%  "'
$ws.Range("D50").Formula = '="0.06218"'
$ws.Range("E50").Formula = '="  +0.41%  "'
$ws.Range("D51").Formula = '="1.853"'
$ws.Range("E51").Formula = '="  +6.11%  "'

$dataRange = $ws.Range("D2:E51")
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

